# "End of EE462 Final Push" -- refresh the measured path-planner metrics
# (pathDist / calcTime / etc.) for the four method blocks on Sheet1 with
# the latest benchmark run's numbers. Only the calcTime (column E) values
# actually moved between runs; methodN / pathDist / numIterations / kGoal /
# kObst / rObst stay the same for every block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (header row 2, data rows 3-5) ---
$ws.Range("E3").Value = 0.14334710000000001
$ws.Range("E4").Value = 0.0174301
$ws.Range("E5").Value = 0.031169100000000002

# --- Block 2 (header row 10, data rows 11-13) ---
$ws.Range("E11").Value = 0.013315799999999999
$ws.Range("E12").Value = 0.0260559
$ws.Range("E13").Value = 0.031995999999999997

# --- Block 3 (header row 18, data rows 19-21) ---
$ws.Range("E19").Value = 0.012299600000000001
$ws.Range("E20").Value = 0.0134176
$ws.Range("E21").Value = 0.027660899999999999

# --- Block 4 (header row 26, data rows 27-29) ---
$ws.Range("E27").Value = 0.0090486000000000004
$ws.Range("E28").Value = 0.0122129
$ws.Range("E29").Value = 0.028326
